# Auto-generated Excel COM-interop script
# Applies the scheduled-runner profit-recalculation update to the Ragnarok_Profits sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 924.5932
$ws.Range("I15").Value = 924.5932
$ws.Range("K15").Value = 2773.7796
$ws.Range("M15").Value = -2604.7796

$ws.Range("H33").Value = 449.85715
$ws.Range("I33").Value = 491.66666
$ws.Range("K33").Value = 491.66666
$ws.Range("M33").Value = -262.66666

$ws.Range("H40").Value = 166668600
$ws.Range("J40").Value = 250001980
$ws.Range("L40").Value = 250001980
$ws.Range("N40").Value = -250002330

$ws.Range("H76").Value = 10909
$ws.Range("I76").Value = 14898
$ws.Range("K76").Value = 14898
$ws.Range("M76").Value = -14583

$ws.Range("H79").Value = 10909
$ws.Range("I79").Value = 14898
$ws.Range("K79").Value = 14898
$ws.Range("M79").Value = -13806

$ws.Range("H135").Value = 4610.125
$ws.Range("I135").Value = 1192.7
$ws.Range("J135").Value = 10305.833
$ws.Range("K135").Value = 10734.3
$ws.Range("L135").Value = 92752.497
$ws.Range("M135").Value = -8199.300000000001
$ws.Range("N135").Value = -97822.497

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8332.209999999999
$ws.Range("I32").Value = 8058.9673
$ws.Range("K32").Value = 8058.9673
$ws.Range("M32").Value = -7771.9673

$ws.Range("H45").Value = 2909.889
$ws.Range("I45").Value = 1507.5385
$ws.Range("J45").Value = 6556
$ws.Range("K45").Value = 1507.5385
$ws.Range("L45").Value = 6556
$ws.Range("M45").Value = -1130.5385
$ws.Range("N45").Value = -7310

$ws.Range("H88").Value = 2736.3845
$ws.Range("J88").Value = 2948.111
$ws.Range("L88").Value = 2948.111
$ws.Range("N88").Value = -3760.111

$ws.Range("H91").Value = 2736.3845
$ws.Range("J91").Value = 2948.111
$ws.Range("L91").Value = 2948.111
$ws.Range("N91").Value = -5756.111

$ws.Range("H122").Value = 2066.8
$ws.Range("J122").Value = 5416.6665
$ws.Range("L122").Value = 16249.9995
$ws.Range("N122").Value = -21149.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3365.139
$ws.Range("I20").Value = 2792.6191
$ws.Range("J20").Value = 4166.6665
$ws.Range("K20").Value = 2792.6191
$ws.Range("L20").Value = 4166.6665
$ws.Range("M20").Value = -2545.6191
$ws.Range("N20").Value = -4660.6665

$ws.Range("H64").Value = 612.4
$ws.Range("I64").Value = 850
$ws.Range("J64").Value = 553
$ws.Range("K64").Value = 850
$ws.Range("L64").Value = 553
$ws.Range("M64").Value = -625
$ws.Range("N64").Value = -1003

$ws.Range("H67").Value = 612.4
$ws.Range("I67").Value = 850
$ws.Range("J67").Value = 553
$ws.Range("K67").Value = 850
$ws.Range("L67").Value = 553
$ws.Range("M67").Value = -70
$ws.Range("N67").Value = -2113

$ws.Range("H81").Value = 19331.666
$ws.Range("J81").Value = 19331.666
$ws.Range("L81").Value = 19331.666
$ws.Range("N81").Value = -21453.666

$ws.Range("H84").Value = 19331.666
$ws.Range("J84").Value = 19331.666
$ws.Range("L84").Value = 57994.99800000001
$ws.Range("N84").Value = -68602.99800000001

$ws.Range("H86").Value = 3409
$ws.Range("I86").Value = 1230.5
$ws.Range("J86").Value = 9944.5
$ws.Range("K86").Value = 1230.5
$ws.Range("L86").Value = 9944.5
$ws.Range("M86").Value = -107.5
$ws.Range("N86").Value = -12190.5

$ws.Range("H89").Value = 3409
$ws.Range("I89").Value = 1230.5
$ws.Range("J89").Value = 9944.5
$ws.Range("K89").Value = 6152.5
$ws.Range("L89").Value = 49722.5
$ws.Range("M89").Value = -536.5
$ws.Range("N89").Value = -60954.5

$ws.Range("H105").Value = 432659.9
$ws.Range("I105").Value = 614757.4399999999
$ws.Range("K105").Value = 614757.4399999999
$ws.Range("M105").Value = -613010.4399999999

$ws.Range("H122").Value = 46081.75
$ws.Range("J122").Value = 45633.332
$ws.Range("L122").Value = 45633.332
$ws.Range("N122").Value = -55433.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34486110
$ws.Range("I31").Value = 71430940
$ws.Range("J31").Value = 4270.933
$ws.Range("K31").Value = 71430940
$ws.Range("L31").Value = 4270.933
$ws.Range("M31").Value = -71430645
$ws.Range("N31").Value = -4860.933

$ws.Range("H34").Value = 34486110
$ws.Range("I34").Value = 71430940
$ws.Range("J34").Value = 4270.933
$ws.Range("K34").Value = 71430940
$ws.Range("L34").Value = 4270.933
$ws.Range("M34").Value = -71430738
$ws.Range("N34").Value = -4674.933

$ws.Range("H58").Value = 2425.682
$ws.Range("I58").Value = 1691.1333
$ws.Range("K58").Value = 1691.1333
$ws.Range("M58").Value = -1488.1333

$ws.Range("H86").Value = 10334.143
$ws.Range("I86").Value = 9900
$ws.Range("J86").Value = 10575.333
$ws.Range("K86").Value = 9900
$ws.Range("L86").Value = 10575.333
$ws.Range("M86").Value = -8777
$ws.Range("N86").Value = -12821.333

$ws.Range("H89").Value = 10334.143
$ws.Range("I89").Value = 9900
$ws.Range("J89").Value = 10575.333
$ws.Range("K89").Value = 49500
$ws.Range("L89").Value = 52876.665
$ws.Range("M89").Value = -43884
$ws.Range("N89").Value = -64108.665

$ws.Range("H105").Value = 5773.125
$ws.Range("I105").Value = 1123.7273
$ws.Range("K105").Value = 1123.7273
$ws.Range("M105").Value = 623.2727

$ws.Range("H122").Value = 2106.56
$ws.Range("I122").Value = 1425.8334
$ws.Range("J122").Value = 3857
$ws.Range("K122").Value = 4277.5002
$ws.Range("L122").Value = 11571
$ws.Range("M122").Value = -1827.5002
$ws.Range("N122").Value = -16471

$ws.Range("H132").Value = 2650.1875
$ws.Range("I132").Value = 2789.6155
$ws.Range("J132").Value = 2046
$ws.Range("K132").Value = 8368.8465
$ws.Range("L132").Value = 6138
$ws.Range("M132").Value = -5838.8465
$ws.Range("N132").Value = -11198

$ws.Range("H136").Value = 2425.682
$ws.Range("I136").Value = 1691.1333
$ws.Range("K136").Value = 5073.3999
$ws.Range("M136").Value = -2523.3999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 10200.4
$ws.Range("I46").Value = 889.6667
$ws.Range("J46").Value = 24166.5
$ws.Range("K46").Value = 2669.0001
$ws.Range("L46").Value = 72499.5
$ws.Range("M46").Value = -2578.0001
$ws.Range("N46").Value = -72681.5

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H113").Value = 878
$ws.Range("J113").Value = 756.5
$ws.Range("L113").Value = 2269.5
$ws.Range("N113").Value = -6609.5

$ws.Range("H117").Value = 3342.8823
$ws.Range("I117").Value = 296.6
$ws.Range("J117").Value = 7694.7144
$ws.Range("K117").Value = 889.8000000000001
$ws.Range("L117").Value = 23084.1432
$ws.Range("M117").Value = 2552.2
$ws.Range("N117").Value = -29968.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5001.385
$ws.Range("I126").Value = 4531.9414
$ws.Range("K126").Value = 13595.8242
$ws.Range("M126").Value = -11125.8242

$ws.Range("H132").Value = 5886891
$ws.Range("I132").Value = 4821.8125
$ws.Range("J132").Value = 100000000
$ws.Range("K132").Value = 14465.4375
$ws.Range("L132").Value = 300000000
$ws.Range("M132").Value = -11935.4375
$ws.Range("N132").Value = -300005060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13214899
$ws.Range("I22").Value = 26424800
$ws.Range("J22").Value = 4999.6
$ws.Range("K22").Value = 26424800
$ws.Range("L22").Value = 4999.6
$ws.Range("M22").Value = -26424505
$ws.Range("N22").Value = -5589.6

$ws.Range("H27").Value = 13214899
$ws.Range("I27").Value = 26424800
$ws.Range("J27").Value = 4999.6
$ws.Range("K27").Value = 26424800
$ws.Range("L27").Value = 4999.6
$ws.Range("M27").Value = -26424693
$ws.Range("N27").Value = -5213.6

$ws.Range("H68").Value = 2195722.2
$ws.Range("I68").Value = 3206955
$ws.Range("K68").Value = 3206955
$ws.Range("M68").Value = -3206206

$ws.Range("H71").Value = 2195722.2
$ws.Range("I71").Value = 3206955
$ws.Range("K71").Value = 16034775
$ws.Range("M71").Value = -16031031

$ws.Range("H122").Value = 4195.5293
$ws.Range("I122").Value = 3383.1628
$ws.Range("K122").Value = 10149.4884
$ws.Range("M122").Value = -7699.4884

$ws.Range("H132").Value = 5353.905
$ws.Range("J132").Value = 6683.357
$ws.Range("L132").Value = 20050.071
$ws.Range("N132").Value = -25110.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 34997.332
$ws.Range("J63").Value = 34997.332
$ws.Range("L63").Value = 34997.332
$ws.Range("N63").Value = -36245.332

$ws.Range("H66").Value = 34997.332
$ws.Range("J66").Value = 34997.332
$ws.Range("L66").Value = 104991.996
$ws.Range("N66").Value = -111231.996

$ws.Range("H81").Value = 1149.1818
$ws.Range("I81").Value = 1064.1
$ws.Range("J81").Value = 2000
$ws.Range("K81").Value = 2128.2
$ws.Range("L81").Value = 4000
$ws.Range("M81").Value = -1067.2
$ws.Range("N81").Value = -6122

$ws.Range("H84").Value = 1149.1818
$ws.Range("I84").Value = 1064.1
$ws.Range("J84").Value = 2000
$ws.Range("K84").Value = 10641
$ws.Range("L84").Value = 20000
$ws.Range("M84").Value = -5337
$ws.Range("N84").Value = -30608

$ws.Range("H100").Value = 1283
$ws.Range("I100").Value = 1349.2727
$ws.Range("J100").Value = 1191.875
$ws.Range("K100").Value = 2698.5454
$ws.Range("L100").Value = 2383.75
$ws.Range("M100").Value = -2157.5454
$ws.Range("N100").Value = -3465.75

$ws.Range("H113").Value = 785.9286
$ws.Range("I113").Value = 595.4286
$ws.Range("J113").Value = 976.4286
$ws.Range("K113").Value = 1786.2858
$ws.Range("L113").Value = 2929.2858
$ws.Range("M113").Value = 383.7142000000001
$ws.Range("N113").Value = -7269.2858

$ws.Range("H132").Value = 1541183.2
$ws.Range("I132").Value = 2974.125
$ws.Range("K132").Value = 8922.375
$ws.Range("M132").Value = -6392.375
